$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '44.193.09'
$ws.Range("E2").Value = '  +1.98%  '
$ws.Range("D3").Value = '2.377.93'
$ws.Range("E3").Value = '  +0.46%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.696'
$ws.Range("E5").Value = '  +7.50%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '243.53'
$ws.Range("E6").Value = '  +4.25%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '76.43'
$ws.Range("E7").Value = '  +6.07%  '
$ws.Range("E8").Value = '  -0.09%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.603'
$ws.Range("E9").Value = '  +25.92%  '
$ws.Range("E10").Value = '  +6.84%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '58.00'
$ws.Range("E11").Value = '  +2.00%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '32.32'
$ws.Range("E12").Value = '  +18.43%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '7.52'
$ws.Range("E13").Value = '  +19.94%  '
$ws.Range("E14").Value = '  +2.73%  '
$ws.Range("D15").Value = '2.733.86'
$ws.Range("E15").Value = '  +0.17%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '17.17'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.928'
$ws.Range("E17").Value = '  +8.07%  '
$ws.Range("D18").Value = '2.385.38'
$ws.Range("E18").Value = '  +0.23%  '
$ws.Range("D19").Value = '44.584.19'
$ws.Range("E19").Value = '  +2.77%  '
$ws.Range("E20").Value = '  +4.58%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.75'
$ws.Range("E21").Value = '  +6.71%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '78.94'
$ws.Range("E22").Value = '  +5.91%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '257.86'
$ws.Range("E23").Value = '  +3.50%  '
$ws.Range("E24").Value = '  +0.05%  '
$ws.Range("E25").Value = '  +4.97%  '
$ws.Range("E26").Value = '  +0.27%  '
$ws.Range("B27").Value = 'ImmutableX'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.81'
$ws.Range("E27").Value = '  +19.74%  '
$ws.Range("B28").Value = 'Cosmos'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.92'
$ws.Range("E28").Value = '  +9.51%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.30'
$ws.Range("E29").Value = '  +1.69%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '23.25'
$ws.Range("E30").Value = '  +3.79%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '175.45'
$ws.Range("E31").Value = '  +0.91%  '
$ws.Range("E32").Value = '  +2.02%  '
$ws.Range("E33").Value = '  +7.47%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.38'
$ws.Range("E34").Value = '  +8.12%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0759'
$ws.Range("E35").Value = '  +9.76%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.35'
$ws.Range("E36").Value = '  +6.24%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.90'
$ws.Range("E37").Value = '  +5.57%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.50'
$ws.Range("E38").Value = '  +2.53%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.61'
$ws.Range("E39").Value = '  +0.23%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0278'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '9.16'
$ws.Range("E41").Value = '  +3.13%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '19.08'
$ws.Range("E42").Value = '  +3.58%  '
$ws.Range("E43").Value = '  +0.04%  '
$ws.Range("E44").Value = '  +17.86%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.22'
$ws.Range("E45").Value = '  +3.40%  '
$ws.Range("B46").Value = 'TrustWalletToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.28'
$ws.Range("E46").Value = '  +5.70%  '
$ws.Range("B47").Value = 'NEARProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.54'
$ws.Range("E47").Value = '  +13.05%  '
$ws.Range("E48").Value = '  +5.95%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '103.14'
$ws.Range("E49").Value = '  +3.21%  '
$ws.Range("E50").Value = '  -0.53%  '
$ws.Range("D51").Value = '1.477.88'
$ws.Range("E51").Value = '  +2.43%  '
